$d = $word.ActiveDocument

# Remove every literal space character from the document body. This is a
# simple "strip spaces" edit: Word collapses the now-empty separator runs,
# merges neighbouring runs that end up sharing identical formatting, and
# drops the spell-check markers / explicit language overrides that no
# longer apply once the text changes.
$range = $d.Content
$range.Find.ClearFormatting()
$range.Find.Replacement.ClearFormatting()
$range.Find.Execute(" ", $true, $false, $false, $false, $false, $true, 1, $false, "", 2)
